$wb = $excel.ActiveWorkbook

# --- Status sheet: Class 010 methods note text change ---
# "At least partially works, but some bugs definitely exist..." ->
# "Buffer to buffer - fine. Between buffers - not so much..."
$wsStatus = $wb.Worksheets.Item("Status")
$ellipsis = [char]0x2026
$wsStatus.Range("B18").Value = "Buffer to buffer - fine. Between buffers - not so much" + $ellipsis

# --- Resolution sheet: fix 1280x1024x16bpp and 1600x1200x16bpp rows ---
# Both go from broken/failing statuses to "WORKING (as much as is implemented)"
# using the same green "WORKING" style already used on row 4 (D4).
$wsRes = $wb.Worksheets.Item("Resolution")

$wsRes.Range("D4").Copy()
$wsRes.Range("D8").PasteSpecial(-4122)
$wsRes.Range("D8").Value = "WORKING (as much as is implemented)"

$wsRes.Range("D4").Copy()
$wsRes.Range("D9").PasteSpecial(-4122)
$wsRes.Range("D9").Value = "WORKING (as much as is implemented)"

# --- Sheet view / selection state updates ---

# Status: scroll/selection moved from A22/B32 to A10/B17
$wsStatus.Range("B17").Select()

# Resolution: now the active/selected tab, selection moved to D9
$wsRes.Activate()
$wsRes.Range("D9").Select()

Write-Output "done"
